$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 915.4545000000001
$ws.Range("I15").Value = 915.4545000000001
$ws.Range("K15").Value = 2746.3635
$ws.Range("M15").Value = -2577.3635
# Row 19
$ws.Range("H19").Value = 1315.4546
$ws.Range("I19").Value = 1367.5
$ws.Range("J19").Value = 1285.7142
$ws.Range("K19").Value = 1367.5
$ws.Range("L19").Value = 1285.7142
$ws.Range("M19").Value = -1192.5
$ws.Range("N19").Value = -1635.7142
# Row 121
$ws.Range("H121").Value = 977.6923
$ws.Range("J121").Value = 892.5
$ws.Range("L121").Value = 2677.5
$ws.Range("N121").Value = -6171.5
# Row 129
$ws.Range("H129").Value = 3294.05
$ws.Range("I129").Value = 6782.6875
$ws.Range("J129").Value = 968.2917
$ws.Range("K129").Value = 20348.0625
$ws.Range("L129").Value = 2904.8751
$ws.Range("M129").Value = -15348.0625
$ws.Range("N129").Value = -12904.8751
# Row 138
$ws.Range("H138").Value = 2292.9727
$ws.Range("I138").Value = 1635.0555
$ws.Range("K138").Value = 4905.166499999999
$ws.Range("M138").Value = 234.8335000000006

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3366.66
$ws.Range("I32").Value = 3090.371
$ws.Range("K32").Value = 3090.371
$ws.Range("M32").Value = -2803.371
# Row 61
$ws.Range("H61").Value = 1452.0571
$ws.Range("I61").Value = 1245.9667
$ws.Range("J61").Value = 2688.6
$ws.Range("K61").Value = 1245.9667
$ws.Range("L61").Value = 2688.6
$ws.Range("M61").Value = -1033.9667
$ws.Range("N61").Value = -3112.6
# Row 80
$ws.Range("H80").Value = 23723.334
$ws.Range("J80").Value = 24061.818
$ws.Range("L80").Value = 24061.818
$ws.Range("N80").Value = -26057.818
# Row 83
$ws.Range("H83").Value = 23723.334
$ws.Range("J83").Value = 24061.818
$ws.Range("L83").Value = 72185.454
$ws.Range("N83").Value = -82169.454
# Row 97
$ws.Range("H97").Value = 26812.77
$ws.Range("I97").Value = 39222.883
$ws.Range("J97").Value = 1992.5385
$ws.Range("K97").Value = 39222.883
$ws.Range("L97").Value = 1992.5385
$ws.Range("M97").Value = -38726.883
$ws.Range("N97").Value = -2984.5385
# Row 110
$ws.Range("H110").Value = 66807410
$ws.Range("I110").Value = 66807410
$ws.Range("K110").Value = 66807410
$ws.Range("M110").Value = -66805365
# Row 132
$ws.Range("H132").Value = 13692.72
$ws.Range("I132").Value = 15935.875
$ws.Range("K132").Value = 47807.625
$ws.Range("M132").Value = -45277.625
# Row 134
$ws.Range("H134").Value = 31070.834
$ws.Range("J134").Value = 31070.834
$ws.Range("L134").Value = 31070.834
$ws.Range("N134").Value = -41210.834
# Row 135
$ws.Range("H135").Value = 14016.765
$ws.Range("I135").Value = 4390
$ws.Range("J135").Value = 15300.333
$ws.Range("K135").Value = 4390
$ws.Range("L135").Value = 15300.333
$ws.Range("M135").Value = 680
$ws.Range("N135").Value = -25440.333
# Row 136
$ws.Range("H136").Value = 1452.0571
$ws.Range("I136").Value = 1245.9667
$ws.Range("J136").Value = 2688.6
$ws.Range("K136").Value = 3737.9001
$ws.Range("L136").Value = 8065.799999999999
$ws.Range("M136").Value = -1187.9001
$ws.Range("N136").Value = -13165.8
# Row 138
$ws.Range("H138").Value = 53996.668
$ws.Range("J138").Value = 53996.668
$ws.Range("L138").Value = 53996.668
$ws.Range("N138").Value = -64276.668
# Row 139
$ws.Range("H139").Value = 59720
$ws.Range("J139").Value = 59720
$ws.Range("L139").Value = 59720
$ws.Range("N139").Value = -70000
# Row 140
$ws.Range("H140").Value = 90160
$ws.Range("J140").Value = 90160
$ws.Range("L140").Value = 90160
$ws.Range("N140").Value = -100520

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
# Row 57
$ws.Range("H57").Value = 54590
$ws.Range("J57").Value = 54590
$ws.Range("L57").Value = 54590
$ws.Range("N57").Value = -56030
# Row 86
$ws.Range("H86").Value = 93801.664
$ws.Range("I86").Value = 124260
$ws.Range("J86").Value = 2426.6667
$ws.Range("K86").Value = 124260
$ws.Range("L86").Value = 2426.6667
$ws.Range("M86").Value = -123137
$ws.Range("N86").Value = -4672.6667
# Row 89
$ws.Range("H89").Value = 93801.664
$ws.Range("I89").Value = 124260
$ws.Range("J89").Value = 2426.6667
$ws.Range("K89").Value = 621300
$ws.Range("L89").Value = 12133.3335
$ws.Range("M89").Value = -615684
$ws.Range("N89").Value = -23365.3335
# Row 94
$ws.Range("H94").Value = 573.48
$ws.Range("I94").Value = 558.4761999999999
$ws.Range("J94").Value = 652.25
$ws.Range("K94").Value = 558.4761999999999
$ws.Range("L94").Value = 652.25
$ws.Range("M94").Value = -107.4761999999999
$ws.Range("N94").Value = -1554.25
# Row 105
$ws.Range("H105").Value = 68402.3
$ws.Range("J105").Value = 101728
$ws.Range("L105").Value = 101728
$ws.Range("N105").Value = -105222
# Row 107
$ws.Range("H107").Value = 71429430
$ws.Range("I107").Value = 100000840
$ws.Range("J107").Value = 925.25
$ws.Range("K107").Value = 100000840
$ws.Range("L107").Value = 925.25
$ws.Range("M107").Value = -99998920
$ws.Range("N107").Value = -4765.25
# Row 132
$ws.Range("H132").Value = 62875
$ws.Range("J132").Value = 62875
$ws.Range("L132").Value = 62875
$ws.Range("N132").Value = -72995
# Row 133
$ws.Range("H133").Value = 134067.5
$ws.Range("J133").Value = 134067.5
$ws.Range("L133").Value = 134067.5
$ws.Range("N133").Value = -144187.5
# Row 134
$ws.Range("H134").Value = 23193.584
$ws.Range("I134").Value = 25337.715
$ws.Range("K134").Value = 76013.145
$ws.Range("M134").Value = -73478.145
# Row 136
$ws.Range("H136").Value = 54590
$ws.Range("J136").Value = 54590
$ws.Range("L136").Value = 54590
$ws.Range("N136").Value = -64790
# Row 137
$ws.Range("H137").Value = 38700
$ws.Range("J137").Value = 38700
$ws.Range("L137").Value = 38700
$ws.Range("N137").Value = -48900
# Row 138
$ws.Range("H138").Value = 63109.168
$ws.Range("J138").Value = 63109.168
$ws.Range("L138").Value = 63109.168
$ws.Range("N138").Value = -73389.16800000001
# Row 139
$ws.Range("H139").Value = 47663.332
$ws.Range("J139").Value = 47663.332
$ws.Range("L139").Value = 47663.332
$ws.Range("N139").Value = -57943.332
# Row 140
$ws.Range("H140").Value = 43175
$ws.Range("J140").Value = 43175
$ws.Range("L140").Value = 43175
$ws.Range("N140").Value = -53535
# Row 141
$ws.Range("H141").Value = 59225
$ws.Range("J141").Value = 59225
$ws.Range("L141").Value = 59225
$ws.Range("N141").Value = -69585

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 949.44446
$ws.Range("I26").Value = 223.75
$ws.Range("J26").Value = 1530
$ws.Range("K26").Value = 671.25
$ws.Range("L26").Value = 4590
$ws.Range("M26").Value = -383.25
$ws.Range("N26").Value = -5166
# Row 86
$ws.Range("H86").Value = 925
$ws.Range("J86").Value = 925
$ws.Range("L86").Value = 2775
$ws.Range("N86").Value = -5147
# Row 89
$ws.Range("H89").Value = 925
$ws.Range("J89").Value = 925
$ws.Range("L89").Value = 8325
$ws.Range("N89").Value = -20181
# Row 131
$ws.Range("H131").Value = 1881.96
$ws.Range("J131").Value = 1955.1158
$ws.Range("L131").Value = 5865.347400000001
$ws.Range("N131").Value = -15945.3474

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 71431360
$ws.Range("I97").Value = 83336024
$ws.Range("J97").Value = 3361
$ws.Range("K97").Value = 83336024
$ws.Range("L97").Value = 3361
$ws.Range("M97").Value = -83335528
$ws.Range("N97").Value = -4353
# Row 132
$ws.Range("H132").Value = 2347.4888
$ws.Range("I132").Value = 1841
$ws.Range("J132").Value = 6399.4
$ws.Range("K132").Value = 5523
$ws.Range("L132").Value = 19198.2
$ws.Range("M132").Value = -2993
$ws.Range("N132").Value = -24258.2
# Row 133
$ws.Range("H133").Value = 53365
$ws.Range("J133").Value = 53365
$ws.Range("L133").Value = 53365
$ws.Range("N133").Value = -63485

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 15585.333
$ws.Range("J20").Value = 15585.333
$ws.Range("L20").Value = 15585.333
$ws.Range("N20").Value = -16037.333
# Row 46
$ws.Range("H46").Value = 533310.5
$ws.Range("I46").Value = 331.66666
$ws.Range("K46").Value = 331.66666
$ws.Range("M46").Value = -143.66666
# Row 64
$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25450
# Row 67
$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26560

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 900.6667
$ws.Range("I107").Value = 900.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2702.0001
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -782.0001000000002
$ws.Range("N107").ClearContents()
# Row 132
$ws.Range("H132").Value = 2620.8298
$ws.Range("I132").Value = 2715.4443
$ws.Range("J132").Value = 2311.182
$ws.Range("K132").Value = 8146.3329
$ws.Range("L132").Value = 6933.545999999999
$ws.Range("M132").Value = -5616.3329
$ws.Range("N132").Value = -11993.546
